$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.268.19'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '3.138.51'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '238.19'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +9.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '636.02'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.95%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.07'
$ws.Range('D7').NumberFormat = "General"
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +3.09%  '
$ws.Range('E8').Value = '  -1.98%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '3.136.64'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.724'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.95%  '
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '36.58'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.88%  '
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').Value = '90.807.33'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '3.721.00'
$ws.Range('D18').Value = '3.137.07'
$ws.Range('E18').Value = '  +1.06%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.74'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.27%  '
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '445.61'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.66'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.67%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.04'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.95'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.61%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '90.22'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.48'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('D28').Value = '3.310.40'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.161'
$ws.Range('D31').NumberFormat = "General"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.95'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +11.67%  '
$ws.Range('B34').Value = 'Stellar'
$ws.Range('C34').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.197'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +23.63%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.83'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '514.39'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('E37').Value = '  +2.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '7.17'
$ws.Range('D38').NumberFormat = "General"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.96%  '
$ws.Range('E39').Value = '  +3.88%  '
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.418'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.12%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.35'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +49.23%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '151.04'
$ws.Range('D47').NumberFormat = "General"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.693'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +9.89%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '45.79'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.54'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +5.62%  '
$ws.Range('E51').Value = '  +3.74%  '
